# Update column G ("K") values on the active worksheet to reflect the
# regenerated save_data (switch from Strike# to K, recalculated std/mean,
# and freshly calculated s_vals written out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 1
    5  = 2
    6  = 2
    7  = 2
    8  = 0
    9  = 2
    10 = 2
    11 = 0
    12 = 2
    13 = 3
    14 = 1
    15 = 0
    16 = 3
    17 = 2
    18 = 0
    19 = 2
    20 = 2
    21 = 1
    22 = 2
    23 = 1
    24 = 3
    25 = 2
    27 = 2
    28 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
